$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the "Periodo Mora" values between rows 16 and 17
$ws.Range("E16").Value = "2008"
$ws.Range("E17").Value = "2009"

# Update "Salario Basico" values for both rows
$ws.Range("G16").Value = 908526
$ws.Range("G17").Value = 908526
